$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue ($ws.Range('D2')) '26.242.00'
Set-TextValue ($ws.Range('E2')) '  -1.78%  '
Set-TextValue ($ws.Range('D3')) '1.583.29'
Set-TextValue ($ws.Range('E3')) '  -1.07%  '
Set-TextValue ($ws.Range('E4')) '  -0.24%  '
Set-TextValue ($ws.Range('D5')) '209.73'
Set-TextValue ($ws.Range('E5')) '  -0.67%  '
Set-TextValue ($ws.Range('D6')) '0.505'
Set-TextValue ($ws.Range('E6')) '  -1.23%  '
Set-TextValue ($ws.Range('E7')) '  -0.23%  '
Set-TextValue ($ws.Range('E8')) '  -1.17%  '
Set-TextValue ($ws.Range('E9')) '  -0.20%  '
Set-TextValue ($ws.Range('D10')) '19.59'
Set-TextValue ($ws.Range('E10')) '  -0.33%  '
Set-TextValue ($ws.Range('E11')) '  +0.36%  '
Set-TextValue ($ws.Range('D12')) '1.806.29'
Set-TextValue ($ws.Range('E12')) '  -1.05%  '
Set-TextValue ($ws.Range('D13')) '1.591.78'
Set-TextValue ($ws.Range('E13')) '  -0.53%  '
Set-TextValue ($ws.Range('D14')) '4.02'
Set-TextValue ($ws.Range('D16')) '64.72'
Set-TextValue ($ws.Range('E16')) '  -0.49%  '
Set-TextValue ($ws.Range('D17')) '26.243.36'
Set-TextValue ($ws.Range('E17')) '  -1.68%  '
Set-TextValue ($ws.Range('D18')) '0.0₃0739'
Set-TextValue ($ws.Range('E18')) '  -0.12%  '
Set-TextValue ($ws.Range('D19')) '7.22'
Set-TextValue ($ws.Range('E19')) '  +0.21%  '
Set-TextValue ($ws.Range('E20')) '  -0.19%  '
Set-TextValue ($ws.Range('D21')) '207.01'
Set-TextValue ($ws.Range('E21')) '  -1.52%  '
Set-TextValue ($ws.Range('E22')) '  -0.86%  '
Set-TextValue ($ws.Range('E23')) '  -3.41%  '
Set-TextValue ($ws.Range('D24')) '8.84'
Set-TextValue ($ws.Range('E24')) '  -1.17%  '
Set-TextValue ($ws.Range('D25')) '144.46'
Set-TextValue ($ws.Range('E25')) '  +0.32%  '
Set-TextValue ($ws.Range('E26')) '  -0.26%  '
Set-TextValue ($ws.Range('D27')) '7.02'
Set-TextValue ($ws.Range('E28')) '  -0.86%  '
Set-TextValue ($ws.Range('D29')) '15.26'
Set-TextValue ($ws.Range('E29')) '  -0.53%  '
Set-TextValue ($ws.Range('E30')) '  -1.50%  '
Set-TextValue ($ws.Range('E31')) '  -1.11%  '
Set-TextValue ($ws.Range('E32')) '  -0.85%  '
Set-TextValue ($ws.Range('E33')) '  -0.71%  '
Set-TextValue ($ws.Range('B34')) 'Maker'
Set-TextValue ($ws.Range('C34')) 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue ($ws.Range('D34')) '1.285.74'
Set-TextValue ($ws.Range('E34')) '  -0.49%  '
Set-TextValue ($ws.Range('B35')) 'WEMIXToken'
Set-TextValue ($ws.Range('C35')) 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue ($ws.Range('D35')) '1.25'
Set-TextValue ($ws.Range('E35')) '  +6.61%  '
Set-TextValue ($ws.Range('E36')) '  -0.09%  '
Set-TextValue ($ws.Range('D37')) '0.610'
Set-TextValue ($ws.Range('E37')) '  +1.40%  '
Set-TextValue ($ws.Range('E38')) '  -1.04%  '
Set-TextValue ($ws.Range('E39')) '  -1.44%  '
Set-TextValue ($ws.Range('D40')) '0.816'
Set-TextValue ($ws.Range('E40')) '  -0.73%  '
Set-TextValue ($ws.Range('D41')) '5.52'
Set-TextValue ($ws.Range('E41')) '  +2.09%  '
Set-TextValue ($ws.Range('D42')) '0.768'
Set-TextValue ($ws.Range('E42')) '  -1.37%  '
Set-TextValue ($ws.Range('E43')) '  -2.75%  '
Set-TextValue ($ws.Range('E44')) '  -1.20%  '
Set-TextValue ($ws.Range('D45')) '1.718.49'
Set-TextValue ($ws.Range('E45')) '  -1.19%  '
Set-TextValue ($ws.Range('D46')) '88.96'
Set-TextValue ($ws.Range('E46')) '  -1.79%  '
Set-TextValue ($ws.Range('D47')) '1.55'
Set-TextValue ($ws.Range('E47')) '  -0.33%  '
Set-TextValue ($ws.Range('E48')) '  +0.72%  '
Set-TextValue ($ws.Range('D49')) '0.0510'
Set-TextValue ($ws.Range('E49')) '  -1.32%  '
Set-TextValue ($ws.Range('E50')) '  -0.28%  '
Set-TextValue ($ws.Range('D51')) '7.42'
Set-TextValue ($ws.Range('E51')) '  +0.20%  '
